$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expense")

# Row 2: Game Controller -> Rent (with trailing space), 70 -> 100, date 45915.22928240741 -> 45931.22928240741
$ws.Range("A2").Value = "Rent "
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 45931.22928240741

# Row 3: Online Shopping -> Movie, 10 -> 8, date 45912.22928240741 -> 45911.22928240741
$ws.Range("A3").Value = "Movie"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 45911.22928240741
